# "Delete All" button now reloads the page: demonstrate the delete/refresh
# cycle by re-populating the contact list rows that reappear after submit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$contacts = @(
    @("Super Corporation", "Alex Reams", "areams@supercorp.com", "630-555-1234"),
    @("Helios Realty",     "Emil Jones", "ejones@email.com",     "645-555-9630")
)

# Rows 2-7: the two contacts repeated three times.
for ($i = 0; $i -lt 6; $i++) {
    $row = $i + 2
    $data = $contacts[$i % 2]
    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
}
